$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 169.63637
$ws.Range("I2").Value = 72.57143000000001
$ws.Range("K2").Value = 72.57143000000001
$ws.Range("M2").Value = 40.42856999999999
$ws.Range("H9").Value = 1671.8889
$ws.Range("I9").Value = 756
$ws.Range("J9").Value = 8999
$ws.Range("K9").Value = 756
$ws.Range("L9").Value = 8999
$ws.Range("M9").Value = -587
$ws.Range("N9").Value = -9337
$ws.Range("H38").Value = 44
$ws.Range("I38").Value = 44
$ws.Range("K38").Value = 132
$ws.Range("M38").Value = 240
$ws.Range("H43").Value = 19005
$ws.Range("J43").Value = 19005
$ws.Range("L43").Value = 19005
$ws.Range("N43").Value = -19143
$ws.Range("H80").Value = 2407.96
$ws.Range("I80").Value = 1574.2667
$ws.Range("J80").Value = 3658.5
$ws.Range("K80").Value = 4722.800099999999
$ws.Range("L80").Value = 10975.5
$ws.Range("M80").Value = -3724.800099999999
$ws.Range("N80").Value = -12971.5
$ws.Range("H83").Value = 2407.96
$ws.Range("I83").Value = 1574.2667
$ws.Range("J83").Value = 3658.5
$ws.Range("K83").Value = 14168.4003
$ws.Range("L83").Value = 32926.5
$ws.Range("M83").Value = -9176.400299999999
$ws.Range("N83").Value = -42910.5
$ws.Range("H101").Value = 1118.8
$ws.Range("I101").Value = 148.5
$ws.Range("J101").Value = 5000
$ws.Range("K101").Value = 445.5
$ws.Range("L101").Value = 15000
$ws.Range("M101").Value = 1176.5
$ws.Range("N101").Value = -18244
$ws.Range("H113").Value = 7682.273
$ws.Range("I113").Value = 5299.8
$ws.Range("J113").Value = 9667.666999999999
$ws.Range("K113").Value = 5299.8
$ws.Range("L113").Value = 9667.666999999999
$ws.Range("M113").Value = -2045.8
$ws.Range("N113").Value = -16175.667
$ws.Range("H121").Value = 2174.5
$ws.Range("J121").Value = 2174.5
$ws.Range("L121").Value = 6523.5
$ws.Range("N121").Value = -10017.5
$ws.Range("H137").Value = 2268.7585
$ws.Range("I137").Value = 590.5
$ws.Range("J137").Value = 3152.0527
$ws.Range("K137").Value = 1771.5
$ws.Range("L137").Value = 9456.158100000001
$ws.Range("M137").Value = 778.5
$ws.Range("N137").Value = -14556.1581
$ws.Range("H138").Value = 2578.6191
$ws.Range("J138").Value = 3208.5
$ws.Range("L138").Value = 9625.5
$ws.Range("N138").Value = -19905.5

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6527
$ws.Range("I61").Value = 4550.6924
$ws.Range("J61").Value = 12950
$ws.Range("K61").Value = 4550.6924
$ws.Range("L61").Value = 12950
$ws.Range("M61").Value = -4338.6924
$ws.Range("N61").Value = -13374
$ws.Range("H98").Value = 5000
$ws.Range("J98").Value = 5000
$ws.Range("L98").Value = 5000
$ws.Range("N98").Value = -10990
$ws.Range("H132").Value = 2716.3914
$ws.Range("I132").Value = 2036.3334
$ws.Range("K132").Value = 6109.0002
$ws.Range("M132").Value = -3579.0002
$ws.Range("H136").Value = 6527
$ws.Range("I136").Value = 4550.6924
$ws.Range("J136").Value = 12950
$ws.Range("K136").Value = 13652.0772
$ws.Range("L136").Value = 38850
$ws.Range("M136").Value = -11102.0772
$ws.Range("N136").Value = -43950

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3113
$ws.Range("I20").Value = 2396.111
$ws.Range("K20").Value = 2396.111
$ws.Range("M20").Value = -2149.111
$ws.Range("H86").Value = 4537.25
$ws.Range("I86").Value = 2838.077
$ws.Range("J86").Value = 7692.857
$ws.Range("K86").Value = 2838.077
$ws.Range("L86").Value = 7692.857
$ws.Range("M86").Value = -1715.077
$ws.Range("N86").Value = -9938.857
$ws.Range("H89").Value = 4537.25
$ws.Range("I89").Value = 2838.077
$ws.Range("J89").Value = 7692.857
$ws.Range("K89").Value = 14190.385
$ws.Range("L89").Value = 38464.285
$ws.Range("M89").Value = -8574.385000000002
$ws.Range("N89").Value = -49696.285
$ws.Range("H94").Value = 1057.95
$ws.Range("I94").Value = 1013.6842
$ws.Range("J94").Value = 1899
$ws.Range("K94").Value = 1013.6842
$ws.Range("L94").Value = 1899
$ws.Range("M94").Value = -562.6842
$ws.Range("N94").Value = -2801
$ws.Range("H105").Value = 8117.3887
$ws.Range("I105").Value = 3511.5
$ws.Range("K105").Value = 3511.5
$ws.Range("M105").Value = -1764.5
$ws.Range("H68").Value = 268
$ws.Range("I68").Value = 268
$ws.Range("K68").Value = 268
$ws.Range("M68").Value = 481
$ws.Range("H71").Value = 268
$ws.Range("I71").Value = 268
$ws.Range("K71").Value = 804
$ws.Range("M71").Value = 2940

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H69").Value = 37750
$ws.Range("I69").Value = 37750
$ws.Range("K69").Value = 37750
$ws.Range("M69").Value = -37001
$ws.Range("H72").Value = 37750
$ws.Range("I72").Value = 37750
$ws.Range("K72").Value = 113250
$ws.Range("M72").Value = -109506
$ws.Range("H132").Value = 5702.2754
$ws.Range("I132").Value = 5522.015
$ws.Range("K132").Value = 16566.045
$ws.Range("M132").Value = -14036.045

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 985.5714
$ws.Range("I23").Value = 70.75
$ws.Range("J23").Value = 1351.5
$ws.Range("K23").Value = 212.25
$ws.Range("L23").Value = 4054.5
$ws.Range("M23").Value = 22.75
$ws.Range("N23").Value = -4524.5
$ws.Range("H37").Value = 161538.23
$ws.Range("J37").Value = 161538.23
$ws.Range("L37").Value = 484614.6900000001
$ws.Range("N37").Value = -484838.6900000001
$ws.Range("H80").Value = 8078.4
$ws.Range("J80").Value = 8124.5
$ws.Range("L80").Value = 24373.5
$ws.Range("N80").Value = -26245.5
$ws.Range("H83").Value = 8078.4
$ws.Range("J83").Value = 8124.5
$ws.Range("L83").Value = 73120.5
$ws.Range("N83").Value = -82480.5
$ws.Range("H97").Value = 690.4
$ws.Range("I97").Value = 249.33333
$ws.Range("J97").Value = 1352
$ws.Range("K97").Value = 747.99999
$ws.Range("L97").Value = 4056
$ws.Range("M97").Value = -251.99999
$ws.Range("N97").Value = -5048
$ws.Range("H98").Value = 999.6667
$ws.Range("I98").Value = 999
$ws.Range("K98").Value = 2997
$ws.Range("M98").Value = -1499
$ws.Range("H140").Value = 4223.467
$ws.Range("I140").Value = 2279.4167
$ws.Range("J140").Value = 11999.667
$ws.Range("K140").Value = 6838.250100000001
$ws.Range("L140").Value = 35999.001
$ws.Range("M140").Value = -1658.250100000001
$ws.Range("N140").Value = -46359.001

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 5500
$ws.Range("I43").Value = 5500
$ws.Range("K43").Value = 5500
$ws.Range("M43").Value = -5349
$ws.Range("H70").Value = 14498.5
$ws.Range("I70").Value = 5545.353
$ws.Range("J70").Value = 36241.855
$ws.Range("K70").Value = 5545.353
$ws.Range("L70").Value = 36241.855
$ws.Range("M70").Value = -5275.353
$ws.Range("N70").Value = -36781.855
$ws.Range("H73").Value = 14498.5
$ws.Range("I73").Value = 5545.353
$ws.Range("J73").Value = 36241.855
$ws.Range("K73").Value = 5545.353
$ws.Range("L73").Value = 36241.855
$ws.Range("M73").Value = -4609.353
$ws.Range("N73").Value = -38113.855
$ws.Range("H93").Value = 25745.7
$ws.Range("J93").Value = 24090.834
$ws.Range("L93").Value = 24090.834
$ws.Range("N93").Value = -27834.834
$ws.Range("H132").Value = 5902.7
$ws.Range("I132").Value = 3126.625
$ws.Range("K132").Value = 9379.875
$ws.Range("M132").Value = -6849.875
$ws.Range("H135").Value = 65799.39999999999
$ws.Range("J135").Value = 65799.39999999999
$ws.Range("L135").Value = 65799.39999999999
$ws.Range("N135").Value = -75939.39999999999

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 900
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H55").Value = 1564598.1
$ws.Range("I55").Value = 2501339.2
$ws.Range("J55").Value = 3362.9167
$ws.Range("K55").Value = 2501339.2
$ws.Range("L55").Value = 3362.9167
$ws.Range("M55").Value = -2501166.2
$ws.Range("N55").Value = -3708.9167
$ws.Range("H61").Value = 3053.6667
$ws.Range("I61").Value = 1004.6429
$ws.Range("K61").Value = 1004.6429
$ws.Range("M61").Value = -802.6429000000001
$ws.Range("H93").Value = 2169.353
$ws.Range("I93").Value = 1741.4
$ws.Range("K93").Value = 1741.4
$ws.Range("M93").Value = -493.4000000000001
$ws.Range("H113").Value = 3053.6667
$ws.Range("I113").Value = 1004.6429
$ws.Range("K113").Value = 1004.6429
$ws.Range("M113").Value = 1165.3571
$ws.Range("H122").Value = 8534.846
$ws.Range("I122").Value = 7094.8
$ws.Range("K122").Value = 21284.4
$ws.Range("M122").Value = -18834.4

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2500
$ws.Range("I62").Value = 2500
$ws.Range("K62").Value = 2500
$ws.Range("M62").Value = -1876
$ws.Range("H65").Value = 2500
$ws.Range("I65").Value = 2500
$ws.Range("K65").Value = 12500
$ws.Range("M65").Value = -9380
$ws.Range("H137").Value = 67997.39999999999
$ws.Range("J137").Value = 67997.39999999999
$ws.Range("L137").Value = 67997.39999999999
$ws.Range("N137").Value = -78197.39999999999
